$d = $word.ActiveDocument

# 1) Merge "Droid Phone: Yes" + "   " (trailing 3 spaces) runs into a single run.
$null = $d.Content.Find.Execute("Droid Phone: Yes   ", $true, $false, $false, $false, $false, $true, 1, $false, "Droid Phone: Yes   ", 2)

# 2) Merge "A final pha" + bookmark("_GoBack") + "se resolves fuzzy spaces..." into a
#    single run and drop the bookmark.
$null = $d.Content.Find.Execute("A final phase resolves fuzzy spaces, and uses x-height normalization to detect lower case letters.", $true, $false, $false, $false, $false, $true, 1, $false, "A final phase resolves fuzzy spaces, and uses x-height normalization to detect lower case letters.", 2)

# 3) Split the lexicon sentence's run after "...by con" and insert a fresh
#    "_GoBack" bookmark (collapsed, zero-length) right before "straining...".
$r = $d.Content
$found = $r.Find.Execute("straining the word by a lexicon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(1)
    $null = $d.Bookmarks.Add("_GoBack", $r)
}
